# Add the "Quelimane" district and its health facilities to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$province    = "Zambézia"
$districtId  = 15
$districtName = "Quelimane"

$facilities = @(
    "CS 17 de Setembro",
    "CS 24 de Julho",
    "CS 4 de Dezembro",
    "CS Chabeco",
    "CS Coalane",
    "CS Incidua",
    "CS Inhangulue",
    "CS Madal",
    "CS Malanha",
    "CS Maquival Rio",
    "CS Maquival Sede",
    "CS Micajune",
    "CS Namuinho",
    "CS Penitenciário",
    "CS Sagariveira",
    "CS Varela",
    "CS Zalala",
    "HG Quelimane"
)

$row = 117
foreach ($facility in $facilities) {
    $ws.Cells.Item($row, 2).Value = $province
    $ws.Cells.Item($row, 3).Value = $districtId
    $ws.Cells.Item($row, 4).Value = $districtName
    $ws.Cells.Item($row, 5).Value = $facility
    $row = $row + 1
}

# Leave the active selection on the last-touched cell, mirroring a manual edit.
$ws.Range("E115").Select()

# Re-registering the filter range creates the (harmless) duplicate
# _xlnm._FilterDatabase_0 defined name that LibreOffice/Excel leave behind
# when the AutoFilter is re-applied after the sheet grows.
$ws.Names.Add("_xlnm._FilterDatabase_0", $ws.Range("B2:E116"))
